$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell {
    param($ws, $row, $col, $text)
    $cell = $ws.Cells.Item($row, $col)
    $origStyle = $cell.Style
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = $origStyle
}

Set-TextCell $ws 2 4 '27.656.36'
Set-TextCell $ws 2 5 '  -0.22%  '
Set-TextCell $ws 3 4 '1.617.89'
Set-TextCell $ws 3 5 '  -0.48%  '
Set-TextCell $ws 4 5 '  -0.28%  '
Set-TextCell $ws 5 4 '209.46'
Set-TextCell $ws 5 5 '  -0.80%  '
Set-TextCell $ws 6 5 '  -1.03%  '
Set-TextCell $ws 7 5 '  -0.21%  '
Set-TextCell $ws 8 4 '23.07'
Set-TextCell $ws 8 5 '  -0.59%  '
Set-TextCell $ws 9 4 '0.255'
Set-TextCell $ws 9 5 '  -1.11%  '
Set-TextCell $ws 10 5 '  -1.03%  '
Set-TextCell $ws 11 5 '  -0.64%  '
Set-TextCell $ws 12 4 '1.846.25'
Set-TextCell $ws 12 5 '  -0.83%  '
Set-TextCell $ws 13 4 '1.627.49'
Set-TextCell $ws 13 5 '  +0.14%  '
Set-TextCell $ws 14 5 '  -1.30%  '
Set-TextCell $ws 15 5 '  -1.04%  '
Set-TextCell $ws 16 4 '64.71'
Set-TextCell $ws 16 5 '  -0.62%  '
Set-TextCell $ws 17 4 '27.665.21'
Set-TextCell $ws 17 5 '  -0.43%  '
Set-TextCell $ws 18 4 '227.38'
Set-TextCell $ws 18 5 '  -1.04%  '
Set-TextCell $ws 19 4 '7.62'
Set-TextCell $ws 19 5 '  +1.58%  '
Set-TextCell $ws 20 5 '  -0.80%  '
Set-TextCell $ws 21 5 '  -0.33%  '
Set-TextCell $ws 22 5 '  -0.87%  '
Set-TextCell $ws 23 5 '  -3.17%  '
Set-TextCell $ws 24 5 '  -1.84%  '
Set-TextCell $ws 25 4 '154.60'
Set-TextCell $ws 25 5 '  +0.13%  '
Set-TextCell $ws 26 4 '6.88'
Set-TextCell $ws 26 5 '  -0.75%  '
Set-TextCell $ws 27 5 '  -0.45%  '
Set-TextCell $ws 28 4 '15.41'
Set-TextCell $ws 28 5 '  -1.10%  '
Set-TextCell $ws 29 5 '  -0.44%  '
Set-TextCell $ws 30 5 '  -1.14%  '
Set-TextCell $ws 31 5 '  -0.47%  '
Set-TextCell $ws 32 5 '  -1.58%  '
Set-TextCell $ws 33 4 '3.07'
Set-TextCell $ws 33 5 '  -0.10%  '
Set-TextCell $ws 34 4 '1.391.62'
Set-TextCell $ws 34 5 '  -0.73%  '
Set-TextCell $ws 35 5 '  +1.39%  '
Set-TextCell $ws 36 5 '  -1.12%  '
Set-TextCell $ws 37 5 '  -1.17%  '
Set-TextCell $ws 38 5 '  +0.54%  '
Set-TextCell $ws 39 5 '  -0.55%  '
Set-TextCell $ws 40 4 '0.841'
Set-TextCell $ws 40 5 '  -2.71%  '
Set-TextCell $ws 41 5 '  -0.42%  '
Set-TextCell $ws 42 5 '  -1.62%  '
Set-TextCell $ws 43 5 '  -0.32%  '
Set-TextCell $ws 44 4 '65.49'
Set-TextCell $ws 44 5 '  -1.45%  '
Set-TextCell $ws 45 4 '5.36'
Set-TextCell $ws 45 5 '  -2.52%  '
Set-TextCell $ws 46 4 '1.756.31'
Set-TextCell $ws 46 5 '  -1.07%  '
Set-TextCell $ws 47 5 '  -7.76%  '
Set-TextCell $ws 48 4 '87.66'
Set-TextCell $ws 48 5 '  -0.09%  '
Set-TextCell $ws 49 5 '  +1.14%  '
Set-TextCell $ws 50 4 '0.0503'
Set-TextCell $ws 50 5 '  -0.72%  '
Set-TextCell $ws 51 4 '7.52'
Set-TextCell $ws 51 5 '  +0.78%  '
